$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45182 (2023-09-13) to 45184 (2023-09-15) for every data row
# (rows 2 through 121).
$firstRow = 2
$lastRow = 121

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
